$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 1539
$ws.Range("K3").Value = 1467
$ws.Range("J4").Value = 1793
$ws.Range("K5").Value = 97
$ws.Range("K6").Value = 1884
$ws.Range("J7").Value = 29261
$ws.Range("K7").Value = 5298

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 100
$ws.Range("K3").Value = 98
$ws.Range("K6").Value = 110
$ws.Range("K7").Value = 332

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 47
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 63
$ws.Range("J4").Value = 56
$ws.Range("K4").Value = 14
$ws.Range("K6").Value = 55
$ws.Range("J7").Value = 1316
$ws.Range("K7").Value = 216

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 26
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 37
$ws.Range("K4").Value = 3
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 38
$ws.Range("K4").Value = 22
$ws.Range("K7").Value = 151
$ws.Range("K8").Value = 332
$ws.Range("K9").Value = 27
$ws.Range("K18").Value = 42
$ws.Range("K19").Value = 135
$ws.Range("K21").Value = 17
$ws.Range("K23").Value = 55
$ws.Range("K29").Value = 241
$ws.Range("K31").Value = 63
$ws.Range("J33").Value = 1316
$ws.Range("K33").Value = 216
$ws.Range("K36").Value = 58
$ws.Range("K37").Value = 178
$ws.Range("K41").Value = 56
$ws.Range("K42").Value = 182
$ws.Range("K43").Value = 52
$ws.Range("K44").Value = 49
$ws.Range("K48").Value = 57
$ws.Range("K51").Value = 64
$ws.Range("K55").Value = 56
$ws.Range("K60").Value = 42
$ws.Range("K63").Value = 17
$ws.Range("K65").Value = 133
$ws.Range("K67").Value = 205
$ws.Range("K76").Value = 72
$ws.Range("K78").Value = 76
$ws.Range("K79").Value = 143
$ws.Range("K83").Value = 107
$ws.Range("K88").Value = 68
$ws.Range("K90").Value = 51
$ws.Range("K91").Value = 52
$ws.Range("K94").Value = 66
$ws.Range("K95").Value = 90
$ws.Range("K96").Value = 72
$ws.Range("J101").Value = 29261
$ws.Range("K101").Value = 5298

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 25
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 63

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 60
$ws.Range("K7").Value = 205

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 63
$ws.Range("K3").Value = 80
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 241

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K5").Value = 8
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 43
$ws.Range("K4").Value = 9
$ws.Range("K6").Value = 80
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 24
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K2").Value = 22
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 52
$ws.Range("K3").Value = 48
$ws.Range("K6").Value = 30
$ws.Range("K7").Value = 143

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 49
$ws.Range("K6").Value = 36
$ws.Range("K7").Value = 151

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 16
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 51

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 16
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 22
